# Updates the cryptocurrency price/volume table (columns D "Price" and
# E "Volume(1h)") on the active sheet with refreshed values scraped by the
# GitHub Actions job. Each cell is written with a leading apostrophe so
# that numeric-looking strings (e.g. "1.00", "598.21") are stored as text
# instead of being auto-coerced to numbers, then ClearFormats() strips the
# implicit "Text" number-format stamp that the apostrophe entry leaves
# behind, keeping the cell's style identical to the untouched cells
# (no explicit style index), matching the original workbook's formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.103.47"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +2.17%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'3.817.25"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +1.29%  "
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'598.21"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +0.33%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'171.60"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  +0.39%  "
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'3.818.86"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'  +1.40%  "
$ws.Range("E7").ClearFormats()
$ws.Range("D9").Value = "'0.525"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  -0.70%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.164"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  -1.24%  "
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'6.54"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  +0.76%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.453"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  -0.76%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'0.0000266"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  -4.11%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'37.04"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  +0.60%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'4.452.90"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  +1.37%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'3.810.85"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +1.21%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'69.064.22"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  +2.14%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'18.34"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  -2.79%  "
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'7.11"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  -2.11%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D21").Value = "'11.19"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  +5.65%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'472.42"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  +0.40%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'0.711"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  -1.72%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'85.09"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  +1.36%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'0.0000146"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  -1.65%  "
$ws.Range("E25").ClearFormats()
$ws.Range("E26").Value = "'  +0.56%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'12.26"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  +0.59%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'10.31"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  -1.34%  "
$ws.Range("E28").ClearFormats()
$ws.Range("E29").Value = "'  +0.13%  "
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = "'3.958.87"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  +1.26%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'2.83"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  -3.09%  "
$ws.Range("E31").ClearFormats()
$ws.Range("D32").Value = "'7.49"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  -3.17%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'2.26"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  +0.05%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'30.40"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  -0.52%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'9.42"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  +2.62%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D37").Value = "'3.767.98"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  +1.11%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'0.103"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  -2.96%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'3.54"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  -8.04%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'0.141"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  +1.74%  "
$ws.Range("E40").ClearFormats()
$ws.Range("E41").Value = "'  +0.54%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'5.88"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  -0.34%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'1.00"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  +0.10%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'0.312"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  -1.38%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D46").Value = "'44.32"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  +13.03%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'1.99"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  +1.32%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'8.66"
$ws.Range("D48").ClearFormats()
$ws.Range("D49").Value = "'46.39"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  +1.05%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'403.70"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  +1.16%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'146.14"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  +2.94%  "
$ws.Range("E51").ClearFormats()
